$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert four new paragraphs before the existing body paragraph:
#    Title, subtitle (parenthetical), authors, blank line.
# ------------------------------------------------------------------
$body = $d.Paragraphs(1)
$body.Range.InsertParagraphBefore()
$body.Range.InsertParagraphBefore()
$body.Range.InsertParagraphBefore()
$body.Range.InsertParagraphBefore()

# Paragraph 1: Title, 18pt (sz/szCs = 36 half-points)
$p1 = $d.Paragraphs(1)
$p1.Range.Text = "Value Proposition"
$p1.Range.Font.Size = 18
$p1.Range.Font.SizeBi = 18

# Paragraph 2: Subtitle in parens, 10pt (sz/szCs = 20 half-points)
$p2 = $d.Paragraphs(2)
$p2.Range.Text = "(AI-enhanced Smart Physical Rehabilitation)"
$p2.Range.Font.Size = 10
$p2.Range.Font.SizeBi = 10

# Paragraph 3: Authors, 10pt
$p3 = $d.Paragraphs(3)
$p3.Range.Text = "Molly Meadows, Noah Rieth, Xian Gao"
$p3.Range.Font.Size = 10
$p3.Range.Font.SizeBi = 10

# Paragraph 4: blank line, 10pt formatting carried in paragraph mark
$p4 = $d.Paragraphs(4)
$p4.Range.Font.Size = 10
$p4.Range.Font.SizeBi = 10

# ------------------------------------------------------------------
# 2. Indent the (former first / now fifth) body paragraph with a
#    first-line indent of 720 twips (0.5").
# ------------------------------------------------------------------
$p5 = $d.Paragraphs(5)
$p5.Format.FirstLineIndent = 36

# ------------------------------------------------------------------
# 3. Text edits inside the body paragraph. Applied right-to-left
#    (rightmost edit first) so each replacement only touches the
#    runs from its own start point onward, leaving the untouched
#    leading runs exactly as they were.
# ------------------------------------------------------------------
function Replace-FirstMatch($oldText, $newText) {
    $full = $d.Content.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        return
    }
    $r = $d.Range($idx, $idx + $oldText.Length)
    $r.Text = $newText
}

# c) "detailed feedback of accuracy and steps" -> "... consisting of an accuracy rating and steps" (rightmost)
Replace-FirstMatch "detailed feedback of accuracy and steps" "detailed feedback consisting of an accuracy rating and steps"

# b) "analyze videos uploaded" -> "analyze uploaded videos" (middle)
Replace-FirstMatch "that will analyze videos uploaded by the client" "that will analyze uploaded videos by the client"

# a) Insert the "physical therapy is expensive" sentence (leftmost)
Replace-FirstMatch "towards total recovery. Our project aims to" "towards total recovery. And since physical therapy is so expensive, any extra time spent with physical therapists procure additional costs. Our project aims to"
